# Auto-generated edit script applying the cryptos price/volume update
# (and the Bittensor/OKB row swap) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.298.58"
$cell.Style = $savedStyle
$cell = $ws.Range("E2")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.29%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D3")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.096.83"
$cell.Style = $savedStyle
$cell = $ws.Range("E3")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.84%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E4")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D5")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "552.08"
$cell.Style = $savedStyle
$cell = $ws.Range("E5")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.90%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D6")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "137.93"
$cell.Style = $savedStyle
$cell = $ws.Range("E6")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -9.05%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E7")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.09%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D8")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.087.76"
$cell.Style = $savedStyle
$cell = $ws.Range("E8")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.91%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E9")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.77%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D10")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.162"
$cell.Style = $savedStyle
$cell = $ws.Range("E10")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.74%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D11")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.54"
$cell.Style = $savedStyle
$cell = $ws.Range("E11")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.60%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D12")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.459"
$cell.Style = $savedStyle
$cell = $ws.Range("E12")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.29%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D13")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "35.11"
$cell.Style = $savedStyle
$cell = $ws.Range("E13")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -6.78%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D14")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000219"
$cell.Style = $savedStyle
$cell = $ws.Range("E14")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.79%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D15")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.602.74"
$cell.Style = $savedStyle
$cell = $ws.Range("E15")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.75%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D16")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.280.75"
$cell.Style = $savedStyle
$cell = $ws.Range("E16")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.16%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E17")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.98%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D18")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.102.71"
$cell.Style = $savedStyle
$cell = $ws.Range("E18")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D19")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "506.24"
$cell.Style = $savedStyle
$cell = $ws.Range("E19")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.71%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D20")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.69"
$cell.Style = $savedStyle
$cell = $ws.Range("E20")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.95%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D21")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.57"
$cell.Style = $savedStyle
$cell = $ws.Range("E21")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -4.22%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E22")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.11%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D23")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.25"
$cell.Style = $savedStyle
$cell = $ws.Range("E23")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.51%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D24")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "77.98"
$cell.Style = $savedStyle
$cell = $ws.Range("E24")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.43%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D25")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.35"
$cell.Style = $savedStyle
$cell = $ws.Range("E25")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -4.84%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E26")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.24%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E27")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.21%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E28")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -7.67%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D29")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $savedStyle
$cell = $ws.Range("E29")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.21%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D30")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.96"
$cell.Style = $savedStyle
$cell = $ws.Range("E30")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -10.02%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D31")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.39"
$cell.Style = $savedStyle
$cell = $ws.Range("E31")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.28%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E32")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -7.63%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E33")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.30%  "
$cell.Style = $savedStyle
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D34")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "57.87"
$cell.Style = $savedStyle
$cell = $ws.Range("E34")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +8.32%  "
$cell.Style = $savedStyle
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D35")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "525.72"
$cell.Style = $savedStyle
$cell = $ws.Range("E35")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -10.64%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D36")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.99"
$cell.Style = $savedStyle
$cell = $ws.Range("E36")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.47%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E37")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -8.04%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E38")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.57%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D39")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.080.74"
$cell.Style = $savedStyle
$cell = $ws.Range("E39")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.29%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D40")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0793"
$cell.Style = $savedStyle
$cell = $ws.Range("E40")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -4.79%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E41")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.74%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D42")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.74"
$cell.Style = $savedStyle
$cell = $ws.Range("E42")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -10.51%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E43")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.36%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D44")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.74"
$cell.Style = $savedStyle
$cell = $ws.Range("E44")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +78.48%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E45")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -2.72%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D47")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "123.03"
$cell.Style = $savedStyle
$cell = $ws.Range("E47")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.89%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E48")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -7.84%  "
$cell.Style = $savedStyle
$cell = $ws.Range("D49")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.34"
$cell.Style = $savedStyle
$cell = $ws.Range("E49")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -7.00%  "
$cell.Style = $savedStyle
$cell = $ws.Range("E50")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -3.14%  "
$cell.Style = $savedStyle
$ws.Range("D51").Value = "0.0₃0508"
$cell = $ws.Range("E51")
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -7.69%  "
$cell.Style = $savedStyle
